$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Add the missing i18n rows (10-15) with the new translation keys.
#    Columns: A = Name (key), B = zh-CN, C = en-US (D = zh-HK, E = fr-FR
#    left blank, same as the existing rows that have no translation yet).
# ---------------------------------------------------------------------------
$rows = @(
    @(10, "navbar.feedback",    "反馈",      "Feedback"),
    @(11, "channel.life",       "生活常用",   "Lifestyle"),
    @(12, "channel.ai",         "人工智能",   "AI"),
    @(13, "channel.image",      "图片视频",   "Media Process"),
    @(14, "channel.developer",  "编程开发",   "Developer"),
    @(15, "channel.external",   "第三方APP",  "External App")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]

    # Match the formatting already used by the existing filled-in rows
    # (2-9): centered text style with the table's fill/border.
    $ws.Range("A9").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)
    $ws.Range("B9").Copy()
    $ws.Range("B$r").PasteSpecial(-4122)
    $ws.Range("C9").Copy()
    $ws.Range("C$r").PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# 2. Column E (fr-FR) never had the table's background fill (only the
#    border). Bring it in line with the rest of the table by applying the
#    same fill color, without touching its existing alignment/border/number
#    format.
# ---------------------------------------------------------------------------
$ws.Range("E1:E26").Interior.ColorIndex = 2
